# Dataset Xa vs ACT: add a header row (Antifactor Xa (IU/mL) / ACT (s)),
# shifting all existing data down by one row, and widen column A to fit
# the new header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 1; existing rows 1-229 (A1:B229)
# move down to A2:B230, and the dimension/used-range grows to A1:B230.
$ws.Range("A1").Insert()

# Populate the new header row with text (stored as shared strings).
$ws.Range("A1").Value = "Antifactor Xa (IU/mL)"
$ws.Range("B1").Value = "ACT (s)"

# Give column A a fixed custom width so the header text is fully visible.
$ws.Columns.Item(1).ColumnWidth = 18.6667

# Reset the view: scroll back to the top and collapse the selection to A1
# (removes the stale topLeftCell="A94" / selection activeCell="E5" state).
$ws.Range("A1").Select()
